# Junction_Flooding_454 refresh: replace the 4 data rows with a new reading
# window, drop the now-unused 6th row, and widen the columns that need the
# extra character of room for the larger values ("custom accuracy" pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New data rows (2-5), 34 columns each (A..AH / Time, J1..J33).
#    These replace the previous rows 2-5 in place; old row 6 is removed below.
# ---------------------------------------------------------------------------
$data = @(
    @(45135.50694444445,13.798,9.137,3.527,29.879,22.444,10.657,31.967,16.87,6.742,10.011,11.733,12.516,3.497,10.903,14.966,9.704000000000001,3.059,1.672,158.575,30.182,10.064,19.331,9.888999999999999,2.87,17.102,8.888999999999999,8.15,9.673999999999999,12.058,3.072,28.934,5.407,12.581),
    @(45135.51388888889,8.055999999999999,5.55,1.481,17.742,13.453,6.218,25.64,9.888999999999999,4.098,5.85,7.031,7.546,2.055,6.391,8.865,5.797,1.38,0.775,89.967,17.97,5.899,11.581,5.976,1.459,12.633,5.211,4.826,5.689,7.212,1.198,23.873,3.143,7.376),
    @(45135.52083333334,5.189,3.609,0.928,11.507,8.673,4,17.765,6.399,2.659,3.731,4.573,4.936,1.331,4.136,5.734,3.796,0.901,0.483,55.627,11.656,3.817,7.497,3.879,0.953,8.465,3.372,3.145,3.694,4.695,0.746,16.48,2.011,4.773),
    @(45135.52777777778,18.17,13.46,1.16,39.73,32.21,14.24,51.89,22.11,9.75,14.4,15.91,16.93,4.59,14.29,20.26,12.12,0.79,0.77,210.12,39.85,13.19,26.69,14.02,2.18,25.79,11.65,10.36,12.19,16.66,0.54,47.08,7.39,16.49)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $i + 2
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# ---------------------------------------------------------------------------
# 2) Drop the old row 6 - the sheet now ends at row 5 (dimension A1:AH5).
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3) Bump column widths by one unit for the columns the diff widened.
#    Excel's ColumnWidth property is offset from the stored <col width>
#    by 5/6 of a character, so "target stored width - 5/6" is what we
#    assign here (7 -> 8 and 8 -> 9 stored widths).
# ---------------------------------------------------------------------------
$widthPlusOne = 7 + (1/6)   # -> stored width 8 (was 7)
$widthPlusOneFromEight = 8 + (1/6)   # -> stored width 9 (was 8)

$colsToEight = @(2, 6, 7, 11, 12, 13, 15, 16, 22, 23, 30, 34)
foreach ($col in $colsToEight) {
    $ws.Columns.Item($col).ColumnWidth = $widthPlusOne
}

$ws.Columns.Item(20).ColumnWidth = $widthPlusOneFromEight
